$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos price list (Price = column D, Volume(1h) = column E).
# Values that Excel would otherwise auto-parse as a plain number (single
# decimal point) are written with a leading apostrophe so they stay text,
# matching the original inline-string cell type.

$ws.Cells.Item(2, 4).Value = '65.444.20'
$ws.Cells.Item(2, 5).Value = '  +2.61%  '
$ws.Cells.Item(3, 4).Value = '3.199.02'
$ws.Cells.Item(3, 5).Value = '  +1.78%  '
$ws.Cells.Item(5, 4).Value = '''599.26'
$ws.Cells.Item(5, 5).Value = '  +2.05%  '
$ws.Cells.Item(6, 4).Value = '''153.92'
$ws.Cells.Item(6, 5).Value = '  +5.79%  '
$ws.Cells.Item(7, 5).Value = '  -0.01%  '
$ws.Cells.Item(8, 4).Value = '3.199.39'
$ws.Cells.Item(8, 5).Value = '  +2.02%  '
$ws.Cells.Item(9, 4).Value = '''0.542'
$ws.Cells.Item(9, 5).Value = '  +2.38%  '
$ws.Cells.Item(10, 5).Value = '  +4.08%  '
$ws.Cells.Item(11, 4).Value = '''6.14'
$ws.Cells.Item(11, 5).Value = '  +6.30%  '
$ws.Cells.Item(12, 4).Value = '''0.472'
$ws.Cells.Item(12, 5).Value = '  +2.90%  '
$ws.Cells.Item(13, 5).Value = '  +3.31%  '
$ws.Cells.Item(14, 4).Value = '''39.36'
$ws.Cells.Item(14, 5).Value = '  +6.53%  '
$ws.Cells.Item(15, 4).Value = '3.732.14'
$ws.Cells.Item(15, 5).Value = '  +1.81%  '
$ws.Cells.Item(16, 5).Value = '  +0.24%  '
$ws.Cells.Item(17, 5).Value = '  +4.59%  '
$ws.Cells.Item(18, 4).Value = '65.124.32'
$ws.Cells.Item(18, 5).Value = '  +2.43%  '
$ws.Cells.Item(19, 4).Value = '3.204.09'
$ws.Cells.Item(19, 5).Value = '  +2.14%  '
$ws.Cells.Item(20, 4).Value = '''485.26'
$ws.Cells.Item(20, 5).Value = '  +4.65%  '
$ws.Cells.Item(21, 5).Value = '  +5.78%  '
$ws.Cells.Item(22, 4).Value = '''0.770'
$ws.Cells.Item(22, 5).Value = '  +5.64%  '
$ws.Cells.Item(23, 4).Value = '''7.93'
$ws.Cells.Item(23, 5).Value = '  +6.51%  '
$ws.Cells.Item(24, 4).Value = '''13.94'
$ws.Cells.Item(24, 5).Value = '  +7.36%  '
$ws.Cells.Item(25, 5).Value = '  +11.71%  '
$ws.Cells.Item(26, 4).Value = '''83.61'
$ws.Cells.Item(26, 5).Value = '  +2.76%  '
$ws.Cells.Item(27, 5).Value = '  +0.32%  '
$ws.Cells.Item(28, 4).Value = '''9.82'
$ws.Cells.Item(28, 5).Value = '  +7.69%  '
$ws.Cells.Item(29, 5).Value = '  +4.00%  '
$ws.Cells.Item(30, 5).Value = '  +3.45%  '
$ws.Cells.Item(31, 5).Value = '  +7.08%  '
$ws.Cells.Item(32, 5).Value = '  -0.04%  '
$ws.Cells.Item(33, 5).Value = '  +9.16%  '
$ws.Cells.Item(34, 4).Value = '''28.59'
$ws.Cells.Item(34, 5).Value = '  +6.10%  '
$ws.Cells.Item(35, 4).Value = '0.0₃0900'
$ws.Cells.Item(35, 5).Value = '  +4.80%  '
$ws.Cells.Item(36, 4).Value = '''3.63'
$ws.Cells.Item(36, 5).Value = '  +7.29%  '
$ws.Cells.Item(37, 5).Value = '  +4.45%  '
$ws.Cells.Item(38, 5).Value = '  +5.78%  '
$ws.Cells.Item(39, 5).Value = '  +3.51%  '
$ws.Cells.Item(40, 4).Value = '''476.52'
$ws.Cells.Item(40, 5).Value = '  +8.05%  '
$ws.Cells.Item(41, 4).Value = '''51.84'
$ws.Cells.Item(41, 5).Value = '  +2.19%  '
$ws.Cells.Item(42, 5).Value = '  +7.69%  '
$ws.Cells.Item(43, 5).Value = '  +9.87%  '
$ws.Cells.Item(44, 4).Value = '''0.0384'
$ws.Cells.Item(44, 5).Value = '  +3.56%  '
$ws.Cells.Item(45, 4).Value = '2.953.97'
$ws.Cells.Item(45, 5).Value = '  +1.69%  '
$ws.Cells.Item(46, 4).Value = '''0.112'
$ws.Cells.Item(46, 5).Value = '  +4.22%  '
$ws.Cells.Item(47, 4).Value = '''38.73'
$ws.Cells.Item(47, 5).Value = '  +5.95%  '
$ws.Cells.Item(48, 4).Value = '''131.78'
$ws.Cells.Item(48, 5).Value = '  +4.80%  '
$ws.Cells.Item(50, 5).Value = '  +5.16%  '
$ws.Cells.Item(51, 5).Value = '  +0.00%  '
